$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet1: remove the duplicate "dh2" / "dh2" / "Invalid Username and Password"
# row (row 7). The shared string "dh2" itself drops out of the shared-string
# table automatically once nothing references it any more.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Rows.Item(7).Delete()

# ---------------------------------------------------------------------------
# Sheet2 ("pythonCode"): widen column A, then append the new sample-code /
# result rows (4-11).
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# Target stored column width is 38.43 characters. Excel's ColumnWidth setter
# round-trips through a pixel grid (MDW=6 for Calibri 11), so feed it the
# character width whose rounded pixel value lands closest to 38.43.
$ws2.Columns.Item(1).ColumnWidth = 37.666666666666664

$searchCode = @"
def search(input_list, num):
if(num in input_list):
print("Element Found")
\b
\b
else:
print("Not Found")
\b
\b
\b
\b
search([12, 23, 45, 67, 6, 90] , 12)
"@

$maxConsecutiveCode = @"
def findMaxConsecutiveOnes(nums) :
count = 0
result = 0
for i in range(0, len(nums)):
if (nums[i] == 0):
count = 0
\b
\b
else:
count+= 1
\b
\b
result = max(result, count)
\b
\b
print(result)
\b
\b
findMaxConsecutiveOnes([1,0,1,1,0,1])
"@

$findNumbersCode = @"
def findNumbers(nums):
c=0
for i in nums:
j=str(i)
x=len(j)
if x%2==0:
c=c+1
\b
\b
\b
\b
print c
return c
findNumbers([12,345,2,6,7896])
"@

$sortedSquaresCode = @"
def sortedSquares(nums):
squares_list = []
for i in range(0, len(nums)):
square = nums[i] * nums[i];
squares_list.append(square)
\b
\b
sorted_squares_list = sorted(squares_list)
print sorted_squares_list;
return sorted_squares_list;
sortedSquares([-7,-3,2,3,11])
"@

$rows = @(
  @{ Row=4;  A=$searchCode;          B="Element Found" },
  @{ Row=5;  A=$searchCode;          B="submission success" },
  @{ Row=6;  A=$maxConsecutiveCode;  B="2" },
  @{ Row=7;  A=$maxConsecutiveCode;  B="submission success" },
  @{ Row=8;  A=$findNumbersCode;     B="2" },
  @{ Row=9;  A=$findNumbersCode;     B="submission success" },
  @{ Row=10; A=$sortedSquaresCode;   B="[4, 9, 9, 49, 121]" },
  @{ Row=11; A=$sortedSquaresCode;   B="submission success" }
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    $cellA = $ws2.Cells.Item($rowNum, 1)
    $cellA.Value = $r.A
    $cellA.Font.Name = "Calibri"
    $cellA.Font.Color = 0
    $cellA.Interior.Color = 16777215
    $cellA.HorizontalAlignment = -4131

    $cellB = $ws2.Cells.Item($rowNum, 2)
    if ($rowNum -ne 4) {
        # Rows 5-11 carry a text-formatted result column (numFmtId 49) - set
        # the format before the value so e.g. "2" is stored as a shared
        # string, not re-interpreted as a number.
        $cellB.NumberFormat = "@"
    }
    $cellB.Value = $r.B

    # Multi-line code in column A makes the host auto-grow the row; AutoFit
    # puts the row back to an implicit (non-customHeight) height, matching
    # the target sheet which has no explicit row heights.
    $ws2.Rows.Item($rowNum).AutoFit()
}

Write-Host "Edit applied"
